$d = $word.ActiveDocument

# --- Candidate Strengths bullets: drop the "(high confidence)." qualifier,
#     add concrete tech examples where indicated ---
$d.Content.Find.Execute(
    "Strong proficiency in Python, Java, or Go (high confidence).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Strong proficiency in Python, Java, or Go", 2) | Out-Null

$d.Content.Find.Execute(
    "Solid understanding of relational databases and NoSQL databases (high confidence).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Solid understanding of relational databases (PostgreSQL, MySQL) and NoSQL databases (MongoDB, Redis)", 2) | Out-Null

$d.Content.Find.Execute(
    "Knowledge of containerization technologies (high confidence).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Knowledge of containerization technologies (Docker, Kubernetes)", 2) | Out-Null

# --- Identified Gaps bullets: rewording ---
$d.Content.Find.Execute(
    "Experience with version control systems (Git) and CI/CD pipelines is only moderately covered.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limited experience with version control systems (Git) and CI/CD pipelines", 2) | Out-Null

# (done via direct Range.Text assignment below, not Find/Replace, so the
#  straight apostrophe in "Bachelor's" is not auto-corrected to a curly one)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $bp = $d.Paragraphs($i)
    if ($bp.Range.Text -eq "Bachelor's degree in a related field is only moderately covered.`r") {
        $bRng = $bp.Range
        $bRng.MoveEnd(1, -1) | Out-Null
        $bRng.Text = "Bachelor's degree not strongly supported"
        break
    }
}

$d.Content.Find.Execute(
    "Knowledge of system design patterns and best practices is the least covered.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Knowledge of system design patterns and best practices is weak", 2) | Out-Null

# --- Insert a new "Risk Flags" section right after that last gap bullet ---
$gapIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Knowledge of system design patterns and best practices is weak`r") {
        $gapIndex = $i
        break
    }
}

$p = $d.Paragraphs($gapIndex)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($gapIndex + 1)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($gapIndex + 2)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs($gapIndex + 3)
$p.Range.InsertParagraphAfter()

# paragraph gapIndex+1 stays as a plain blank separator paragraph
$d.Paragraphs($gapIndex + 1).Style = "Normal"

# paragraph gapIndex+2 becomes the "Risk Flags" heading
$riskHeading = $d.Paragraphs($gapIndex + 2)
$riskHeading.Range.Text = "Risk Flags"
$riskHeading.Style = "Heading1"

# paragraph gapIndex+3 becomes the first risk bullet
$riskBullet1 = $d.Paragraphs($gapIndex + 3)
$riskBullet1.Range.Text = "⚠️ Low confidence in version control and CI/CD experience"
$riskBullet1.Style = "ListBullet"

# paragraph gapIndex+4 becomes the second risk bullet
$riskBullet2 = $d.Paragraphs($gapIndex + 4)
$riskBullet2.Range.Text = "⚠️ Weak educational background support"
$riskBullet2.Style = "ListBullet"

# --- Drafted Communication body: drop the two QA-suggestion lines ---
$lineBreak = [char]11
$oldBody = "- Clarify the candidate's experience with version control systems and CI/CD pipelines." + $lineBreak + "- Verify the educational background to ensure it meets the job requirements." + $lineBreak + $lineBreak
$d.Content.Find.Execute($oldBody, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Next Steps / QA Suggestions bullets: rewording ---
$d.Content.Find.Execute(
    "Clarify the candidate's experience with version control systems and CI/CD pipelines.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Consider a technical interview to assess practical skills in version control and CI/CD", 2) | Out-Null

$d.Content.Find.Execute(
    "Verify the educational background to ensure it meets the job requirements.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Verify educational qualifications and consider additional certifications", 2) | Out-Null

Write-Host "Done"
